$wb = $excel.ActiveWorkbook

# --- Add the new "Renames" sheet as the last tab ---------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$renames = $wb.Worksheets.Add($null, $lastSheet)
$renames.Name = "Renames"

# Header row.
$renames.Range("A1").Value = "old name"
$renames.Range("B1").Value = "new name"

# old name -> new name pairs (mirrors the renamed components in the
# Composition sheet: "Cylinder" -> "Big Cylinder", "Lug Nut" -> "Locking Nut").
$renames.Range("A2").Value = "Cylinder"
$renames.Range("B2").Value = "Big Cylinder"

$renames.Range("A3").Value = "Lug Nut"
$renames.Range("B3").Value = "Locking Nut"

# Selection on the new sheet lands just below the data, on A4.
$renames.Range("A4").Select()

# --- Scroll the Composition sheet's view down one row ----------------------
# (topLeftCell becomes A2, selection stays E11). Must be done while the
# Composition sheet is active, so switch to it, scroll, then switch back.
$composition = $wb.Worksheets.Item("Composition")
$composition.Activate()
$excel.ActiveWindow.ScrollRow = 2

# --- Make "Renames" the active/visible tab again ----------------------------
$renames.Activate()
